$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "a"
$ws.Range("C8").Value = "b"
$ws.Range("D7").Value = "c"
$ws.Range("D9").Value = "d"

$ws.Range("G9").Select()
